$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- C11 / C12: add a solid-yellow highlight fill ----------------------
# (done first so the "yellow" cellXf slot is allocated before the
# "theme" one, matching the style-table order of the target workbook)
$ws.Range("C11").Interior.Color = 65535
$ws.Range("C12").Interior.Color = 65535

# --- C4 / C9: 0 -> 1e-9, highlight fill = theme Accent4, lighter 60% --
# Setting a plain color first, then overwriting with ThemeColor on the
# same cell, keeps the interop layer from leaving a spare placeholder
# fill behind (a quirk of how it resolves Interior.ThemeColor).
$ws.Range("C4").Value = 0.0000000010000000000000001
$ws.Range("C4").Interior.Color = 65535
$ws.Range("C4").Interior.ThemeColor = 8

$ws.Range("C9").Value = 0.0000000010000000000000001
$ws.Range("C9").Interior.Color = 65535
$ws.Range("C9").Interior.ThemeColor = 8

# --- New row 14: "Water" / "water" = -1e-9, same theme highlight ------
$ws.Range("A14").Value2 = "Water"
$ws.Range("B14").Value2 = "water"
$ws.Range("C14").Value = -0.0000000010000000000000001
$ws.Range("A14:C14").Interior.Color = 65535
$ws.Range("A14:C14").Interior.ThemeColor = 8

# --- Selection / active cell as last saved ------------------------------
$ws.Range("B18").Select()
